$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.166.76'
$ws.Range('E2').Value = '  -0.36%  '

$ws.Range('D3').Value = '2.270.63'
$ws.Range('E3').Value = '  +1.05%  '

$ws.Range('E4').Value = '  -0.78%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '498.94'
$ws.Range('E5').Value = '  +0.80%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.63'
$ws.Range('E6').Value = '  +0.99%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.20%  '

$ws.Range('E8').Value = '  -0.73%  '

$ws.Range('E9').Value = '  -0.11%  '

$ws.Range('E10').Value = '  +0.57%  '

$ws.Range('E11').Value = '  +3.29%  '

$ws.Range('E12').Value = '  +1.61%  '

$ws.Range('D13').Value = '2.670.68'
$ws.Range('E13').Value = '  +0.23%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.62'
$ws.Range('E14').Value = '  +4.25%  '

$ws.Range('D15').Value = '54.122.83'
$ws.Range('E15').Value = '  -0.27%  '

$ws.Range('E16').Value = '  +0.26%  '

$ws.Range('D17').Value = '2.276.90'
$ws.Range('E17').Value = '  -0.53%  '

$ws.Range('E18').Value = '  +2.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.13'
$ws.Range('E19').Value = '  +1.72%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '302.76'
$ws.Range('E20').Value = '  -0.61%  '

$ws.Range('E21').Value = '  -2.51%  '

$ws.Range('E22').Value = '  +0.62%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.08'
$ws.Range('E23').Value = '  -2.98%  '

$ws.Range('E24').Value = '  -0.28%  '

$ws.Range('E25').Value = '  -1.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.27'
$ws.Range('E26').Value = '  +2.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.37'
$ws.Range('E27').Value = '  -0.23%  '

$ws.Range('E28').Value = '  +0.59%  '

$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.90'
$ws.Range('E29').Value = '  +0.49%  '

$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0683'
$ws.Range('E30').Value = '  -0.56%  '

$ws.Range('E31').Value = '  +0.97%  '

$ws.Range('E32').Value = '  -0.03%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.74'
$ws.Range('E33').Value = '  +0.78%  '

$ws.Range('E34').Value = '  +10.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  +0.47%  '

$ws.Range('E36').Value = '  -1.63%  '

$ws.Range('E37').Value = '  +0.84%  '

$ws.Range('E38').Value = '  -0.99%  '

$ws.Range('E39').Value = '  -0.11%  '

$ws.Range('E40').Value = '  +0.49%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.80'
$ws.Range('E41').Value = '  -1.61%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '124.72'
$ws.Range('E42').Value = '  -3.65%  '

$ws.Range('E43').Value = '  +1.99%  '

$ws.Range('E44').Value = '  -0.40%  '

$ws.Range('E45').Value = '  -1.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '238.29'
$ws.Range('E46').Value = '  -1.76%  '

$ws.Range('E47').Value = '  -0.83%  '

$ws.Range('E48').Value = '  +0.76%  '

$ws.Range('E49').Value = '  +0.49%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.16'
$ws.Range('E50').Value = '  -1.09%  '

$ws.Range('E51').Value = '  -0.34%  '
